# "Data using new search query"
# Refresh the publication-year / count table with results from the new
# Web of Science search. The new query returns a longer history (back to
# 1969 instead of 2015) and higher counts for every year already present.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Year -> Count, most recent first (row 2 = 2025 ... row 37 = 1969).
$rows = @(
    @("2025", 27),
    @("2024", 1374),
    @("2023", 1312),
    @("2022", 1250),
    @("2021", 1142),
    @("2020", 1029),
    @("2019", 842),
    @("2018", 679),
    @("2017", 571),
    @("2016", 485),
    @("2015", 371),
    @("2014", 297),
    @("2013", 278),
    @("2012", 186),
    @("2011", 143),
    @("2010", 114),
    @("2009", 75),
    @("2008", 73),
    @("2007", 62),
    @("2006", 33),
    @("2005", 19),
    @("2004", 18),
    @("2003", 16),
    @("2002", 14),
    @("2001", 13),
    @("2000", 9),
    @("1999", 9),
    @("1998", 4),
    @("1997", 5),
    @("1996", 5),
    @("1995", 3),
    @("1994", 6),
    @("1993", 2),
    @("1991", 2),
    @("1979", 1),
    @("1969", 1)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $year = $rows[$i][0]
    $count = $rows[$i][1]

    # Leading apostrophe forces the year to be stored as text, matching the
    # original "Publication Years" column (years are labels, not numbers).
    $ws.Cells.Item($r, 1).Value = "'" + $year
    $ws.Cells.Item($r, 2).Value = $count
}
